$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.827.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.67%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.620.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.92%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'594.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.34%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'149.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.28%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.71%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +4.21%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.27%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.150"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.18%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'27.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.58%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'3.092.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.91%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'63.711.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.70%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0000148"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.14%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.668.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.07%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.16%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'348.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.69%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.59%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.64%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'66.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.51%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +12.97%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -0.13%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'9.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.89%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "'Kaspa"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'0.163"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.09%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Aptos"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'8.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.50%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'542.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.08%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +0.05%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.52%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.0₃0846"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.99%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -0.33%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.85%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'168.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.54%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'FirstDigitalUSD"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.07%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'PolygonEcosystemToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.23%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.69%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'19.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.53%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'168.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.54%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'39.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.45%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +4.00%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.0588"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.82%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'21.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.25%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.48%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +11.91%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.37%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +0.38%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'19.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.83%  "
$ws.Range("E51").Style = "Normal"
